$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col4a5"
$ws.Range("C2").Value = "Cd93"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1320423333333333
$ws.Range("H2").Value = 0.396127
$ws.Range("I2").Value = 0.01362486282338958
$ws.Range("J2").Value = 0.01362486282338958
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 122.328922
$ws.Range("N2").Value = 366.986766
$ws.Range("O2").Value = 0.9783373008518612
$ws.Range("P2").Value = 0.9783373008518613
$ws.Range("Q2").Value = 16.15259629503133
$ws.Range("R2").Value = 145.373366655282
$ws.Range("S2").Value = 0.01332971151911183
$ws.Range("T2").Value = 0.01332971151911183
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col4a5"
$ws.Range("C3").Value = "Cd93"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1320423333333333
$ws.Range("H3").Value = 0.396127
$ws.Range("I3").Value = 0.01362486282338958
$ws.Range("J3").Value = 0.01362486282338958
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3863573333333334
$ws.Range("N3").Value = 1.159072
$ws.Range("O3").Value = 0.003089929874945324
$ws.Range("P3").Value = 0.003089929874945324
$ws.Range("Q3").Value = 0.05101552379377779
$ws.Range("R3").Value = 0.4591397141440001
$ws.Range("S3").Value = 0.00004209987068002335
$ws.Range("T3").Value = 0.00004209987068002335
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col4a5"
$ws.Range("C4").Value = "Cd93"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1320423333333333
$ws.Range("H4").Value = 0.396127
$ws.Range("I4").Value = 0.01362486282338958
$ws.Range("J4").Value = 0.01362486282338958
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.322294
$ws.Range("N4").Value = 6.966882000000001
$ws.Range("O4").Value = 0.0185727692731934
$ws.Range("P4").Value = 0.0185727692731934
$ws.Range("Q4").Value = 0.3066411184460001
$ws.Range("R4").Value = 2.759770066014001
$ws.Range("S4").Value = 0.0002530514335977251
$ws.Range("T4").Value = 0.0002530514335977251
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col4a5"
$ws.Range("C5").Value = "Cd93"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.730541000000001
$ws.Range("H5").Value = 23.191623
$ws.Range("I5").Value = 0.7976802440297347
$ws.Range("J5").Value = 0.7976802440297348
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 122.328922
$ws.Range("N5").Value = 366.986766
$ws.Range("O5").Value = 0.9783373008518612
$ws.Range("P5").Value = 0.9783373008518613
$ws.Range("Q5").Value = 945.6687470068022
$ws.Range("R5").Value = 8511.018723061219
$ws.Range("S5").Value = 0.7804003368869046
$ws.Range("T5").Value = 0.7804003368869048
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col4a5"
$ws.Range("C6").Value = "Cd93"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7.730541000000001
$ws.Range("H6").Value = 23.191623
$ws.Range("I6").Value = 0.7976802440297347
$ws.Range("J6").Value = 0.7976802440297348
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3863573333333334
$ws.Range("N6").Value = 1.159072
$ws.Range("O6").Value = 0.003089929874945324
$ws.Range("P6").Value = 0.003089929874945324
$ws.Range("Q6").Value = 2.986751205984001
$ws.Range("R6").Value = 26.88076085385601
$ws.Range("S6").Value = 0.002464776016681154
$ws.Range("T6").Value = 0.002464776016681154
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col4a5"
$ws.Range("C7").Value = "Cd93"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7.730541000000001
$ws.Range("H7").Value = 23.191623
$ws.Range("I7").Value = 0.7976802440297347
$ws.Range("J7").Value = 0.7976802440297348
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.322294
$ws.Range("N7").Value = 6.966882000000001
$ws.Range("O7").Value = 0.0185727692731934
$ws.Range("P7").Value = 0.0185727692731934
$ws.Range("Q7").Value = 17.95258898105401
$ws.Range("R7").Value = 161.573300829486
$ws.Range("S7").Value = 0.01481513112614887
$ws.Range("T7").Value = 0.01481513112614887
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col4a5"
$ws.Range("C8").Value = "Cd93"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.828694666666667
$ws.Range("H8").Value = 5.486084
$ws.Range("I8").Value = 0.1886948931468756
$ws.Range("J8").Value = 0.1886948931468756
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 122.328922
$ws.Range("N8").Value = 366.986766
$ws.Range("O8").Value = 0.9783373008518612
$ws.Range("P8").Value = 0.9783373008518613
$ws.Range("Q8").Value = 223.7022472404826
$ws.Range("R8").Value = 2013.320225164344
$ws.Range("S8").Value = 0.1846072524458446
$ws.Range("T8").Value = 0.1846072524458447
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col4a5"
$ws.Range("C9").Value = "Cd93"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.828694666666667
$ws.Range("H9").Value = 5.486084
$ws.Range("I9").Value = 0.1886948931468756
$ws.Range("J9").Value = 0.1886948931468756
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3863573333333334
$ws.Range("N9").Value = 1.159072
$ws.Range("O9").Value = 0.003089929874945324
$ws.Range("P9").Value = 0.003089929874945324
$ws.Range("Q9").Value = 0.7065295948942223
$ws.Range("R9").Value = 6.358766354048001
$ws.Range("S9").Value = 0.0005830539875841466
$ws.Range("T9").Value = 0.0005830539875841466
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col4a5"
$ws.Range("C10").Value = "Cd93"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.828694666666667
$ws.Range("H10").Value = 5.486084
$ws.Range("I10").Value = 0.1886948931468756
$ws.Range("J10").Value = 0.1886948931468756
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.322294
$ws.Range("N10").Value = 6.966882000000001
$ws.Range("O10").Value = 0.0185727692731934
$ws.Range("P10").Value = 0.0185727692731934
$ws.Range("Q10").Value = 4.246766652232
$ws.Range("R10").Value = 38.220899870088
$ws.Range("S10").Value = 0.003504586713446804
$ws.Range("T10").Value = 0.003504586713446804
